$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 72157
$ws.Range("J3").Value = 72157
$ws.Range("L3").Value = 72157
$ws.Range("N3").Value = -72385
$ws.Range("H6").Value = 470.9565
$ws.Range("I6").Value = 470.9565
$ws.Range("K6").Value = 1412.8695
$ws.Range("M6").Value = -1300.8695
$ws.Range("H12").Value = 1548.1666
$ws.Range("J12").Value = 9000
$ws.Range("L12").Value = 9000
$ws.Range("N12").Value = -9340
$ws.Range("H18").Value = 1538
$ws.Range("I18").Value = 1422.5
$ws.Range("J18").Value = 2000
$ws.Range("K18").Value = 1422.5
$ws.Range("L18").Value = 2000
$ws.Range("M18").Value = -1138.5
$ws.Range("N18").Value = -2568
$ws.Range("H21").Value = 10749.5
$ws.Range("I21").Value = 1000
$ws.Range("J21").Value = 13999.333
$ws.Range("K21").Value = 1000
$ws.Range("L21").Value = 13999.333
$ws.Range("M21").Value = -532
$ws.Range("N21").Value = -14935.333
$ws.Range("H23").Value = 10749.5
$ws.Range("I23").Value = 1000
$ws.Range("J23").Value = 13999.333
$ws.Range("K23").Value = 1000
$ws.Range("L23").Value = 13999.333
$ws.Range("M23").Value = -766
$ws.Range("N23").Value = -14467.333
$ws.Range("H31").Value = 150
$ws.Range("I31").Value = 150
$ws.Range("K31").Value = 450
$ws.Range("M31").Value = -220
$ws.Range("H32").Value = 20005000
$ws.Range("J32").Value = 25005750
$ws.Range("L32").Value = 25005750
$ws.Range("N32").Value = -25006402
$ws.Range("H38").Value = 1001.2308
$ws.Range("I38").Value = 1001.2308
$ws.Range("K38").Value = 3003.6924
$ws.Range("M38").Value = -2631.6924
$ws.Range("H52").Value = 3443.1667
$ws.Range("J52").Value = 3350
$ws.Range("L52").Value = 10050
$ws.Range("N52").Value = -10370
$ws.Range("H64").Value = 7834.6
$ws.Range("J64").Value = 9306.5
$ws.Range("L64").Value = 9306.5
$ws.Range("N64").Value = -9802.5
$ws.Range("H67").Value = 7834.6
$ws.Range("J67").Value = 9306.5
$ws.Range("L67").Value = 9306.5
$ws.Range("N67").Value = -11022.5
$ws.Range("H87").Value = 28809.523
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 28809.523
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 28809.523
$ws.Range("M87").ClearContents()
$ws.Range("N87").Value = -31305.523
$ws.Range("H90").Value = 28809.523
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 28809.523
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 86428.569
$ws.Range("M90").ClearContents()
$ws.Range("N90").Value = -98908.569
$ws.Range("H92").Value = 462.8846
$ws.Range("J92").Value = 428.33334
$ws.Range("L92").Value = 428.33334
$ws.Range("N92").Value = -2924.33334
$ws.Range("H100").Value = 8473.352999999999
$ws.Range("I100").Value = 8264.166999999999
$ws.Range("J100").Value = 8587.454
$ws.Range("K100").Value = 8264.166999999999
$ws.Range("L100").Value = 8587.454
$ws.Range("M100").Value = -7723.166999999999
$ws.Range("N100").Value = -9669.454
$ws.Range("H102").Value = 72157
$ws.Range("J102").Value = 72157
$ws.Range("L102").Value = 72157
$ws.Range("N102").Value = -78647
$ws.Range("H106").Value = 648.25
$ws.Range("I106").Value = 648.25
$ws.Range("K106").Value = 648.25
$ws.Range("M106").Value = -17.25
$ws.Range("H111").Value = 2935.8
$ws.Range("I111").Value = 3064.3333
$ws.Range("J111").Value = 2743
$ws.Range("K111").Value = 9192.999899999999
$ws.Range("L111").Value = 8229
$ws.Range("M111").Value = -6125.999899999999
$ws.Range("N111").Value = -14363
$ws.Range("H116").Value = 8183.769
$ws.Range("I116").Value = 4749.5
$ws.Range("K116").Value = 4749.5
$ws.Range("M116").Value = -1307.5
$ws.Range("H125").Value = 1442.125
$ws.Range("I125").Value = 709.5
$ws.Range("J125").Value = 2174.75
$ws.Range("K125").Value = 6385.5
$ws.Range("L125").Value = 19572.75
$ws.Range("M125").Value = -3925.5
$ws.Range("N125").Value = -24492.75
$ws.Range("H128").Value = 39818.184
$ws.Range("J128").Value = 39818.184
$ws.Range("L128").Value = 39818.184
$ws.Range("N128").Value = -49778.184
$ws.Range("H132").Value = 12827.682
$ws.Range("J132").Value = 166890.67
$ws.Range("L132").Value = 500672.01
$ws.Range("N132").Value = -505732.01
$ws.Range("H135").Value = 11369492
$ws.Range("I135").Value = 17244796
$ws.Range("J135").Value = 10571.267
$ws.Range("K135").Value = 155203164
$ws.Range("L135").Value = 95141.40299999999
$ws.Range("M135").Value = -155200629
$ws.Range("N135").Value = -100211.403
$ws.Range("H138").Value = 3286.3547
$ws.Range("I138").Value = 1954.2941
$ws.Range("J138").Value = 4903.857
$ws.Range("K138").Value = 5862.8823
$ws.Range("L138").Value = 14711.571
$ws.Range("M138").Value = -722.8823000000002
$ws.Range("N138").Value = -24991.571
$ws.Range("H141").Value = 5758.1064
$ws.Range("I141").Value = 4038.8333
$ws.Range("J141").Value = 20200
$ws.Range("K141").Value = 12116.4999
$ws.Range("L141").Value = 60600
$ws.Range("M141").Value = -6936.499899999999
$ws.Range("N141").Value = -70960

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 4163.1665
$ws.Range("I45").Value = 3467.9092
$ws.Range("J45").Value = 6075.125
$ws.Range("K45").Value = 3467.9092
$ws.Range("L45").Value = 6075.125
$ws.Range("M45").Value = -3090.9092
$ws.Range("N45").Value = -6829.125
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("M60").ClearContents()
$ws.Range("N60").ClearContents()
$ws.Range("H74").Value = 3849.5293
$ws.Range("J74").Value = 4138.2856
$ws.Range("L74").Value = 4138.2856
$ws.Range("N74").Value = -5886.2856
$ws.Range("H76").Value = 50000
$ws.Range("J76").Value = 50000
$ws.Range("L76").Value = 50000
$ws.Range("N76").Value = -50676
$ws.Range("H77").Value = 3849.5293
$ws.Range("J77").Value = 4138.2856
$ws.Range("L77").Value = 20691.428
$ws.Range("N77").Value = -29427.428
$ws.Range("H79").Value = 50000
$ws.Range("J79").Value = 50000
$ws.Range("L79").Value = 50000
$ws.Range("N79").Value = -52340
$ws.Range("H82").Value = 44999
$ws.Range("J82").Value = 44999
$ws.Range("L82").Value = 44999
$ws.Range("N82").Value = -45721
$ws.Range("H85").Value = 44999
$ws.Range("J85").Value = 44999
$ws.Range("L85").Value = 44999
$ws.Range("N85").Value = -47495
$ws.Range("H88").Value = 1760.4445
$ws.Range("I88").Value = 1748.75
$ws.Range("J88").Value = 1763.7858
$ws.Range("K88").Value = 1748.75
$ws.Range("L88").Value = 1763.7858
$ws.Range("M88").Value = -1342.75
$ws.Range("N88").Value = -2575.7858
$ws.Range("H91").Value = 1760.4445
$ws.Range("I91").Value = 1748.75
$ws.Range("J91").Value = 1763.7858
$ws.Range("K91").Value = 1748.75
$ws.Range("L91").Value = 1763.7858
$ws.Range("M91").Value = -344.75
$ws.Range("N91").Value = -4571.7858
$ws.Range("H97").Value = 5231.5386
$ws.Range("I97").Value = 4000.8333
$ws.Range("K97").Value = 4000.8333
$ws.Range("M97").Value = -3504.8333
$ws.Range("H110").Value = 3904.7437
$ws.Range("I110").Value = 3971.4
$ws.Range("K110").Value = 3971.4
$ws.Range("M110").Value = -1926.4
$ws.Range("H122").Value = 6089.2144
$ws.Range("I122").Value = 5892.857
$ws.Range("J122").Value = 6285.5713
$ws.Range("K122").Value = 17678.571
$ws.Range("L122").Value = 18856.7139
$ws.Range("M122").Value = -15228.571
$ws.Range("N122").Value = -23756.7139

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H37").Value = 998
$ws.Range("I37").Value = 998
$ws.Range("K37").Value = 998
$ws.Range("M37").Value = -861
$ws.Range("H76").Value = 35000
$ws.Range("J76").Value = 35000
$ws.Range("L76").Value = 35000
$ws.Range("N76").Value = -35630
$ws.Range("H79").Value = 35000
$ws.Range("J79").Value = 35000
$ws.Range("L79").Value = 35000
$ws.Range("N79").Value = -37184
$ws.Range("H86").Value = 1663.6
$ws.Range("I86").Value = 1003.6667
$ws.Range("K86").Value = 1003.6667
$ws.Range("M86").Value = 119.3333
$ws.Range("H89").Value = 1663.6
$ws.Range("I89").Value = 1003.6667
$ws.Range("K89").Value = 5018.3335
$ws.Range("M89").Value = 597.6665000000003
$ws.Range("H94").Value = 4240.8
$ws.Range("I94").Value = 3824.75
$ws.Range("K94").Value = 3824.75
$ws.Range("M94").Value = -3373.75
$ws.Range("H105").Value = 4787.2856
$ws.Range("I105").Value = 3912.2
$ws.Range("K105").Value = 3912.2
$ws.Range("M105").Value = -2165.2
$ws.Range("H107").Value = 3711.375
$ws.Range("I107").Value = 3058.111
$ws.Range("J107").Value = 5671.1665
$ws.Range("K107").Value = 3058.111
$ws.Range("L107").Value = 5671.1665
$ws.Range("M107").Value = -1138.111
$ws.Range("N107").Value = -9511.166499999999
$ws.Range("H134").Value = 7089.7188
$ws.Range("I134").Value = 2191.2693
$ws.Range("J134").Value = 28316.334
$ws.Range("K134").Value = 6573.8079
$ws.Range("L134").Value = 84949.00199999999
$ws.Range("M134").Value = -4038.8079
$ws.Range("N134").Value = -90019.00199999999

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 8912.5
$ws.Range("I16").Value = 10595
$ws.Range("J16").Value = 500
$ws.Range("K16").Value = 10595
$ws.Range("L16").Value = 500
$ws.Range("M16").Value = -10308
$ws.Range("N16").Value = -1074
$ws.Range("H31").Value = 6115.4165
$ws.Range("I31").Value = 9359.333000000001
$ws.Range("J31").Value = 2871.5
$ws.Range("K31").Value = 9359.333000000001
$ws.Range("L31").Value = 2871.5
$ws.Range("M31").Value = -9064.333000000001
$ws.Range("N31").Value = -3461.5
$ws.Range("H34").Value = 6115.4165
$ws.Range("I34").Value = 9359.333000000001
$ws.Range("J34").Value = 2871.5
$ws.Range("K34").Value = 9359.333000000001
$ws.Range("L34").Value = 2871.5
$ws.Range("M34").Value = -9157.333000000001
$ws.Range("N34").Value = -3275.5
$ws.Range("H70").Value = 25000
$ws.Range("J70").Value = 25000
$ws.Range("L70").Value = 25000
$ws.Range("N70").Value = -25630
$ws.Range("H73").Value = 25000
$ws.Range("J73").Value = 25000
$ws.Range("L73").Value = 25000
$ws.Range("N73").Value = -27184
$ws.Range("H97").Value = 44179
$ws.Range("J97").Value = 0
$ws.Range("L97").Value = 0
$ws.Range("N97").ClearContents()
$ws.Range("H113").Value = 8912.5
$ws.Range("I113").Value = 10595
$ws.Range("J113").Value = 500
$ws.Range("K113").Value = 10595
$ws.Range("L113").Value = 500
$ws.Range("M113").Value = -8425
$ws.Range("N113").Value = -4840
$ws.Range("H141").Value = 52500
$ws.Range("J141").Value = 52500
$ws.Range("L141").Value = 52500
$ws.Range("N141").Value = -62860

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 798.1905
$ws.Range("J5").Value = 489.13333
$ws.Range("L5").Value = 1467.39999
$ws.Range("N5").Value = -1691.39999
$ws.Range("H22").Value = 333333340
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("H27").Value = 333333340
$ws.Range("I27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("H36").Value = 793.6
$ws.Range("J36").Value = 1751.5
$ws.Range("L36").Value = 5254.5
$ws.Range("N36").Value = -5592.5
$ws.Range("H64").Value = 4726
$ws.Range("J64").Value = 5499.5
$ws.Range("L64").Value = 16498.5
$ws.Range("N64").Value = -17038.5
$ws.Range("H67").Value = 4726
$ws.Range("J67").Value = 5499.5
$ws.Range("L67").Value = 16498.5
$ws.Range("N67").Value = -18370.5
$ws.Range("H75").Value = 15000
$ws.Range("J75").Value = 15000
$ws.Range("L75").Value = 45000
$ws.Range("N75").Value = -46996
$ws.Range("H78").Value = 15000
$ws.Range("J78").Value = 15000
$ws.Range("L78").Value = 135000
$ws.Range("N78").Value = -144984
$ws.Range("H101").Value = 3999
$ws.Range("J101").Value = 3999
$ws.Range("L101").Value = 11997
$ws.Range("N101").Value = -16865
$ws.Range("H107").Value = 1404.3226
$ws.Range("J107").Value = 2362.9092
$ws.Range("L107").Value = 7088.7276
$ws.Range("N107").Value = -10928.7276
$ws.Range("H117").Value = 1202
$ws.Range("I117").Value = 494.5
$ws.Range("J117").Value = 1359.2222
$ws.Range("K117").Value = 1483.5
$ws.Range("L117").Value = 4077.6666
$ws.Range("M117").Value = 1958.5
$ws.Range("N117").Value = -10961.6666
$ws.Range("H135").Value = 798.1905
$ws.Range("J135").Value = 489.13333
$ws.Range("L135").Value = 4402.19997
$ws.Range("N135").Value = -9472.19997
$ws.Range("H139").Value = 76926660
$ws.Range("I139").Value = 100003060
$ws.Range("J139").Value = 5333
$ws.Range("K139").Value = 300009180
$ws.Range("L139").Value = 15999
$ws.Range("M139").Value = -300004040
$ws.Range("N139").Value = -26279

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H9").Value = 2741.4
$ws.Range("J9").Value = 1733.3334
$ws.Range("L9").Value = 1733.3334
$ws.Range("N9").Value = -2073.3334
$ws.Range("H97").Value = 10021.4375
$ws.Range("I97").Value = 2638.7
$ws.Range("K97").Value = 2638.7
$ws.Range("M97").Value = -2142.7
$ws.Range("H126").Value = 2879.64
$ws.Range("I126").Value = 2313.625
$ws.Range("J126").Value = 3885.889
$ws.Range("K126").Value = 6940.875
$ws.Range("L126").Value = 11657.667
$ws.Range("M126").Value = -4470.875
$ws.Range("N126").Value = -16597.667
$ws.Range("H132").Value = 6083.1665
$ws.Range("I132").Value = 5812.871
$ws.Range("J132").Value = 7759
$ws.Range("K132").Value = 17438.613
$ws.Range("L132").Value = 23277
$ws.Range("M132").Value = -14908.613
$ws.Range("N132").Value = -28337

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 2500.5
$ws.Range("I16").Value = 2500.5
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2500.5
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2330.5
$ws.Range("N16").ClearContents()
$ws.Range("H22").Value = 3393.7273
$ws.Range("J22").Value = 3500
$ws.Range("L22").Value = 3500
$ws.Range("N22").Value = -4090
$ws.Range("H27").Value = 3393.7273
$ws.Range("J27").Value = 3500
$ws.Range("L27").Value = 3500
$ws.Range("N27").Value = -3714
$ws.Range("H40").Value = 5174.3335
$ws.Range("I40").Value = 5517.857
$ws.Range("K40").Value = 5517.857
$ws.Range("M40").Value = -5381.857
$ws.Range("H46").Value = 4751.7036
$ws.Range("I46").Value = 3633.3333
$ws.Range("J46").Value = 4891.5
$ws.Range("K46").Value = 3633.3333
$ws.Range("L46").Value = 4891.5
$ws.Range("M46").Value = -3445.3333
$ws.Range("N46").Value = -5267.5
$ws.Range("H55").Value = 309.5909
$ws.Range("I55").Value = 441.53845
$ws.Range("J55").Value = 119
$ws.Range("K55").Value = 441.53845
$ws.Range("L55").Value = 119
$ws.Range("M55").Value = -268.53845
$ws.Range("N55").Value = -465
$ws.Range("H68").Value = 3278
$ws.Range("J68").Value = 3278
$ws.Range("L68").Value = 3278
$ws.Range("N68").Value = -4776
$ws.Range("H71").Value = 3278
$ws.Range("J71").Value = 3278
$ws.Range("L71").Value = 16390
$ws.Range("N71").Value = -23878
$ws.Range("H82").Value = 3234.6667
$ws.Range("I82").Value = 3065
$ws.Range("K82").Value = 3065
$ws.Range("M82").Value = -2704
$ws.Range("H85").Value = 3234.6667
$ws.Range("I85").Value = 3065
$ws.Range("K85").Value = 3065
$ws.Range("M85").Value = -1817
$ws.Range("H93").Value = 2484.4443
$ws.Range("I93").Value = 2337.1428
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 2337.1428
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -1089.1428
$ws.Range("N93").Value = -5496
$ws.Range("H100").Value = 95190.07000000001
$ws.Range("I100").Value = 111405.55
$ws.Range("J100").Value = 35733.332
$ws.Range("K100").Value = 111405.55
$ws.Range("L100").Value = 35733.332
$ws.Range("M100").Value = -110864.55
$ws.Range("N100").Value = -36815.332
$ws.Range("H132").Value = 3841.6428
$ws.Range("I132").Value = 3912.5386
$ws.Range("J132").Value = 2920
$ws.Range("K132").Value = 11737.6158
$ws.Range("L132").Value = 8760
$ws.Range("M132").Value = -9207.6158
$ws.Range("N132").Value = -13820
$ws.Range("H136").Value = 3580.3076
$ws.Range("I136").Value = 4549.3335
$ws.Range("J136").Value = 1400
$ws.Range("K136").Value = 13648.0005
$ws.Range("L136").Value = 4200
$ws.Range("M136").Value = -11098.0005
$ws.Range("N136").Value = -9300

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 45000
$ws.Range("I70").Value = 25000
$ws.Range("J70").Value = 55000
$ws.Range("K70").Value = 25000
$ws.Range("L70").Value = 55000
$ws.Range("M70").Value = -24685
$ws.Range("N70").Value = -55630
$ws.Range("H73").Value = 45000
$ws.Range("I73").Value = 25000
$ws.Range("J73").Value = 55000
$ws.Range("K73").Value = 25000
$ws.Range("L73").Value = 55000
$ws.Range("M73").Value = -23908
$ws.Range("N73").Value = -57184
$ws.Range("H81").Value = 18510.896
$ws.Range("I81").Value = 18826.691
$ws.Range("J81").Value = 18254.312
$ws.Range("K81").Value = 37653.382
$ws.Range("L81").Value = 36508.624
$ws.Range("M81").Value = -36592.382
$ws.Range("N81").Value = -38630.624
$ws.Range("H84").Value = 18510.896
$ws.Range("I84").Value = 18826.691
$ws.Range("J84").Value = 18254.312
$ws.Range("K84").Value = 188266.91
$ws.Range("L84").Value = 182543.12
$ws.Range("M84").Value = -182962.91
$ws.Range("N84").Value = -193151.12
$ws.Range("H96").Value = 46664.832
$ws.Range("I96").Value = 87413
$ws.Range("K96").Value = 87413
$ws.Range("M96").Value = -86040
$ws.Range("H107").Value = 656
$ws.Range("I107").Value = 703.5333000000001
$ws.Range("K107").Value = 2110.5999
$ws.Range("M107").Value = -190.5999000000002
$ws.Range("H136").Value = 3778.718
$ws.Range("I136").Value = 2352
$ws.Range("K136").Value = 7056
$ws.Range("M136").Value = -4506
